# Updates the genetic-algorithm iteration log (rows 2-21) with the
# latest run's results: best chromosome (C), max/min/avg fitness (D/E/F).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C holds 32-bit chromosome strings (e.g. "1111...0001"); Excel
# would otherwise read these digit strings as numbers, so the range is
# formatted as Text before the values are written.
$ws.Range("C2:C21").NumberFormat = "@"

$ws.Range("C2").Value = "111101000010110000110111000001"
$ws.Range("D2").Value = 0.9097338043760219
$ws.Range("E2").Value = 0.06909138166141518
$ws.Range("F2").Value = 0.3985843274486325
$ws.Range("C3").Value = "111101011001110110010100111101"
$ws.Range("D3").Value = 0.9205169757907861
$ws.Range("E3").Value = 0.4628926776578759
$ws.Range("F3").Value = 0.7773979506887685
$ws.Range("C4").Value = "111101011001110110010100111101"
$ws.Range("D4").Value = 0.9205169757907861
$ws.Range("E4").Value = 0.4953339706513745
$ws.Range("F4").Value = 0.8704501706201248
$ws.Range("C5").Value = "111101011010110000110111000001"
$ws.Range("D5").Value = 0.9209454765071264
$ws.Range("E5").Value = 0.9093079202541157
$ws.Range("F5").Value = 0.9130195768543501
$ws.Range("C6").Value = "111101011010110000010100111101"
$ws.Range("D6").Value = 0.9209415797519598
$ws.Range("E6").Value = 0.9092681414866814
$ws.Range("F6").Value = 0.9118521012717249
$ws.Range("C7").Value = "111101011010110000010100111101"
$ws.Range("D7").Value = 0.9209415797519598
$ws.Range("E7").Value = 0.9092681414866814
$ws.Range("F7").Value = 0.9118480999705818
$ws.Range("C8").Value = "111101011010110000010100100001"
$ws.Range("D8").Value = 0.9209415297019471
$ws.Range("E8").Value = 0.9092681414866814
$ws.Range("F8").Value = 0.9129606498761735
$ws.Range("C9").Value = "111101011001110110110111000001"
$ws.Range("D9").Value = 0.9205208716475425
$ws.Range("E9").Value = 0.9092681414866814
$ws.Range("F9").Value = 0.910718603932749
$ws.Range("C10").Value = "111101000010110000110111001101"
$ws.Range("D10").Value = 0.909733825695106
$ws.Range("E10").Value = 0.9092681414866814
$ws.Range("F10").Value = 0.9095933309166673
$ws.Range("C11").Value = "111101000010110000110111000001"
$ws.Range("D11").Value = 0.9097338043760219
$ws.Range("E11").Value = 0.9092681414866814
$ws.Range("F11").Value = 0.909593328784759
$ws.Range("C12").Value = "111101000011110000110111000001"
$ws.Range("D12").Value = 0.9101995864746523
$ws.Range("E12").Value = 0.9088025978066306
$ws.Range("F12").Value = 0.909593350494713
$ws.Range("C13").Value = "111101000011110000110111000001"
$ws.Range("D13").Value = 0.9101995864746523
$ws.Range("E13").Value = 0.9088025978066306
$ws.Range("F13").Value = 0.9097334507106265
$ws.Range("C14").Value = "111101000011110000110111000001"
$ws.Range("D14").Value = 0.9101995864746523
$ws.Range("E14").Value = 0.9088025978066306
$ws.Range("F14").Value = 0.909686483072581
$ws.Range("C15").Value = "111101000011110000110111000001"
$ws.Range("D15").Value = 0.9101995864746523
$ws.Range("E15").Value = 0.9088025978066306
$ws.Range("F15").Value = 0.9097322823263987
$ws.Range("C16").Value = "111101000011110000110111000001"
$ws.Range("D16").Value = 0.9101995864746523
$ws.Range("E16").Value = 0.9097299100941983
$ws.Range("F16").Value = 0.9098250136482366
$ws.Range("C17").Value = "111101000011110000110111000001"
$ws.Range("D17").Value = 0.9101995864746523
$ws.Range("E17").Value = 0.9097299100941983
$ws.Range("F17").Value = 0.9097780460167908
$ws.Range("C18").Value = "111101000011110000110111000001"
$ws.Range("D18").Value = 0.9101995864746523
$ws.Range("E18").Value = 0.9097299100941983
$ws.Range("F18").Value = 0.9098246242266539
$ws.Range("C19").Value = "111101000011110000110111000001"
$ws.Range("D19").Value = 0.9101995864746523
$ws.Range("E19").Value = 0.9097299100941983
$ws.Range("F19").Value = 0.9097776565886084
$ws.Range("C20").Value = "111101000011110000110111000001"
$ws.Range("D20").Value = 0.9101995864746523
$ws.Range("E20").Value = 0.9097299100941983
$ws.Range("F20").Value = 0.9098715918646993
$ws.Range("C21").Value = "111101000011110000110111000001"
$ws.Range("D21").Value = 0.9101995864746523
$ws.Range("E21").Value = 0.9097299100941983
$ws.Range("F21").Value = 0.9098715918646991
